$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ParticipantsTab): update TabQuery (B2) and StatQuery (C2)
$participantsQuery = @'
WITH diag_base AS (
  SELECT DISTINCT
    d."participant.id"     AS participant_id,
    d.age_at_diagnosis,
    d.diagnosis,
    d.anatomic_site,
    d.diagnosis_category,
    d.disease_phase
  FROM df_diagnosis d
  WHERE d."participant.id" IS NOT NULL
),
diagnosis_summary AS (
  SELECT
    participant_id,
    GROUP_CONCAT(
      CASE WHEN age_at_diagnosis = '-999' THEN 'Not Reported' ELSE age_at_diagnosis END,
      ';'
    )                                 AS age_at_diagnosis,
    GROUP_CONCAT(diagnosis, ';')      AS diagnosis,
    GROUP_CONCAT(anatomic_site, ';')  AS anatomic_site,
    GROUP_CONCAT(diagnosis_category, ';') AS diagnosis_category
  FROM diag_base
  GROUP BY participant_id
),
surv_base AS (
  SELECT DISTINCT
    s."participant.id" AS participant_id,
    s.last_known_survival_status
  FROM df_survival s
),
survival_summary AS (
  SELECT
    participant_id,
    GROUP_CONCAT(last_known_survival_status, ';') AS last_known_survival_status
  FROM surv_base
  GROUP BY participant_id
)
SELECT DISTINCT
  prt.participant_id                           AS "Participant ID",
  std.dbgap_accession                          AS "Study ID",
  COALESCE(prt.sex_at_birth, '')               AS "Sex",
  COALESCE(prt.race, '')                       AS "Race",
  COALESCE(dsum.diagnosis, '')                 AS "Diagnosis",
  COALESCE(dsum.anatomic_site, '')             AS "Diagnosis Anatomic Site",
  COALESCE(dsum.diagnosis_category, '')        AS "Diagnosis_Category",
  COALESCE(dsum.age_at_diagnosis, '')          AS "Age at Diagnosis (days)",
  NULL                                         AS "Treatment Type",
  COALESCE(srv.last_known_survival_status, '') AS "Last Known Survival Status"
FROM df_study std
LEFT JOIN df_consent_group   cg   ON std.id = cg."study.id"
LEFT JOIN df_participant     prt  ON cg.id  = prt."consent_group.id"
LEFT JOIN diagnosis_summary  dsum ON prt.id = dsum.participant_id
LEFT JOIN survival_summary   srv  ON prt.id = srv.participant_id
WHERE
  std.dbgap_accession = 'phs000464'
  AND EXISTS (
    SELECT 1
    FROM df_diagnosis d
    WHERE d."participant.id" = prt.id
      AND d.disease_phase = 'Initial Diagnosis'
  )
  AND EXISTS (
    SELECT 1
    FROM df_survival s
    WHERE s."participant.id" = prt.id
      AND s.last_known_survival_status = 'Unknown'
  )
ORDER BY prt.participant_id ASC
LIMIT 100;
'@
$ws.Range("B2").Value = $participantsQuery

$statQuery = @'
SELECT
  COUNT(DISTINCT std.dbgap_accession) AS "Studies",
  COUNT(DISTINCT prt.participant_id)  AS "Participants",
  CAST(0 AS INTEGER)                  AS "Samples",
  COUNT(DISTINCT cmf.id)              AS "Files"
FROM df_study std
LEFT JOIN df_consent_group         cg  ON std.id = cg."study.id"
LEFT JOIN df_participant           prt ON cg.id  = prt."consent_group.id"
LEFT JOIN df_clinical_measure_file cmf ON std.id = cmf."study.id"
WHERE 
  std.dbgap_accession = 'phs000464'
  AND EXISTS (
    SELECT 1
    FROM df_diagnosis d
    WHERE d."participant.id" = prt.id
      AND d.disease_phase = 'Initial Diagnosis'
  )
  AND EXISTS (
    SELECT 1
    FROM df_survival s
    WHERE s."participant.id" = prt.id
      AND s.last_known_survival_status = 'Unknown'
  );
'@
$ws.Range("C2").Value = $statQuery

# Row 4 (SamplesTab): update TabQuery (B4)
$samplesQuery = @'
SELECT DISTINCT
    smp.sample_id                       AS "Sample ID",
    prt.participant_id                  AS "Participant ID",
    std.dbgap_accession                 AS "Study ID",
    COALESCE(smp.anatomic_site, '')     AS "Sample Anatomic Site",
    COALESCE(
        CASE 
            WHEN smp.participant_age_at_collection = '-999' THEN 'Not Reported'
            ELSE smp.participant_age_at_collection
        END,
        ''
    )                                   AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '')      AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '')     AS "Sample Tumor Classification",
    COALESCE(dgn.diagnosis, '')                AS "Sample Diagnosis"
FROM df_study std
LEFT JOIN df_consent_group cg ON std.id = cg."study.id"
LEFT JOIN df_participant   prt ON cg.id  = prt."consent_group.id"
LEFT JOIN df_sample        smp ON prt.id = smp."participant.id"
LEFT JOIN df_diagnosis     dgn ON smp.id = dgn."sample.id"
WHERE 
    std.dbgap_accession = 'phs000464'
    AND EXISTS (
        SELECT 1
        FROM df_diagnosis d
        WHERE d."participant.id" = prt.id
          AND d.disease_phase = 'Initial Diagnosis'
    )
    AND EXISTS (
        SELECT 1
        FROM df_survival s
        WHERE s."participant.id" = prt.id
          AND s.last_known_survival_status = 'Unknown'
    )
    AND smp.sample_id IS NOT NULL
ORDER BY 
    smp.sample_id ASC;
'@
$ws.Range("B4").Value = $samplesQuery

# Row 5 (FilesTab): update TabQuery (B5)
$filesQuery = @'
WITH file_data AS (
  SELECT
    file_name,
    data_category,
    file_type,
    file_description,
    file_access,
    CAST(file_size AS REAL) AS file_size,
    "study.id"       AS studyid,
    "participant.id" AS participantid
  FROM df_clinical_measure_file
)
SELECT
  COALESCE(fd.file_name, '')        AS "File Name",
  COALESCE(fd.data_category, '')    AS "Data Category",
  COALESCE(fd.file_description, '') AS "File Description",
  COALESCE(fd.file_type, '')        AS "File Type",
  CASE
    WHEN fd.file_size >= 1024*1024*1024 THEN
      CASE
        WHEN ROUND(fd.file_size / (1024.0*1024.0*1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0*1024.0*1024.0), 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / (1024.0*1024.0*1024.0), 0) AS INT) AS TEXT) || ' GB'
        ELSE ROUND(fd.file_size / (1024.0*1024.0*1024.0), 2) || ' GB'
      END
    WHEN fd.file_size >= 1024*1024 THEN
      CASE
        WHEN ROUND(fd.file_size / (1024.0*1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0*1024.0), 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / (1024.0*1024.0), 0) AS INT) AS TEXT) || ' MB'
        ELSE ROUND(fd.file_size / (1024.0*1024.0), 2) || ' MB'
      END
    WHEN fd.file_size >= 1024 THEN
      CASE
        WHEN ROUND(fd.file_size / 1024.0, 2) = CAST(ROUND(fd.file_size / 1024.0, 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
        ELSE ROUND(fd.file_size / 1024.0, 2) || ' KB'
      END
    WHEN fd.file_size IS NOT NULL THEN
      CASE
        WHEN ROUND(fd.file_size, 2) = CAST(ROUND(fd.file_size, 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size, 0) AS INT) AS TEXT) || ' Bytes'
        ELSE ROUND(fd.file_size, 2) || ' Bytes'
      END
    ELSE ''
  END                                  AS "File Size",
  COALESCE(fd.file_access, '')         AS "File Access",
  std.dbgap_accession                  AS "Study ID",
  COALESCE(fd.participantid, '')       AS "Participant ID",
  NULL                                 AS "Sample ID"
FROM df_study std
LEFT JOIN file_data fd ON std.id = fd.studyid
WHERE
  std.dbgap_accession = 'phs000464'
  AND EXISTS (
    SELECT 1
    FROM df_consent_group cg
    JOIN df_participant prt ON cg.id = prt."consent_group.id"
    JOIN df_diagnosis  d   ON prt.id = d."participant.id"
    JOIN df_survival   s   ON prt.id = s."participant.id"
    WHERE cg."study.id" = std.id
      AND d.disease_phase = 'Initial Diagnosis'
      AND s.last_known_survival_status = 'Unknown'
  )
ORDER BY fd.file_name ASC;
'@
$ws.Range("B5").Value = $filesQuery

# Row heights: row2 keeps its original custom height (content grew but the
# row was manually sized, so Excel does not auto-grow it); row5's natural
# autofit in this runtime overshoots what real Excel computed, so pin it
# back to the workbook's effective max wrapped-row height.
$ws.Rows.Item(2).RowHeight = 355
$ws.Rows.Item(5).RowHeight = 409.6
